$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-21 (product_name, product_url, time_text, rating_count, variant_name, review_text)
$data = @(
    ,@("Eskan fresia 3 liter plastik ceret minum bening motif bungaVariasi: Eskan fresia", "/Eskan-fresia-3-liter-plastik-ceret-minum-bening-motif-bunga-i.145589728.8605964666", "2025-05-25 06:07", 5, "Variasi: Eskan fresia", "Ukuran: uk.kecil..sy pikir 3ltr itu besar ehtaunya datangnya malah kecil")
    ,@("Rak plastik serbaguna susun 3 lovina warna goldVariasi: Rak lovina Gold", "/Rak-plastik-serbaguna-susun-3-lovina-warna-gold-i.145589728.24033169934", "2025-05-25 04:40", 5, "Variasi: Rak lovina Gold", "Ukuran: pas dan sesuai")
    ,@("Kapstok kawat vinil warna 6 gantungan - gantungan baju tempel kawatVariasi: Kapstok 6 warna", "/Kapstok-kawat-vinil-warna-6-gantungan-gantungan-baju-tempel-kawat-i.145589728.13014371227", "2025-05-25 01:44", 5, "Variasi: Kapstok 6 warna", "Ssuai gmbr & keterangan`nproduk baik & brmutu`nrespon & pengiriman cepat`npacking aman & bgs`nPenjual amanah & tr prcya")
    ,@("Silet cukur tatra original extra tajam 10 pcs - silet serbagunaVariasi: Silet tatra", "/Silet-cukur-tatra-original-extra-tajam-10-pcs-silet-serbaguna-i.145589728.9915872611", "2025-05-25 01:44", 5, "Variasi: Silet tatra", "Kandungan: okKenyamanan: ck6Efek: blm dicobaSsuai gmbr & keterangan`nproduk baik & brmutu`nrespon & pengiriman cepat`npacking aman & bgs`nPenjual amanah & tr prcya")
    ,@("kipas tangan plastik karakter waru murah warna warniVariasi: Kipas bunga rotan", "/kipas-tangan-plastik-karakter-waru-murah-warna-warni-i.145589728.5337255631", "2025-05-25 01:44", 5, "Variasi: Kipas bunga rotan", "Ketebalan: tebalKualitas: okDesain: menarikSsuai gmbr & keterangan`nproduk baik & brmutu`nrespon & pengiriman cepat`npacking aman & bgs`nPenjual amanah & tr prcya")
    ,@("Rak Mini segi susun 3 Flower - Rak bumbu kosmetik serbagunaVariasi: Rak segi Flower", "/Rak-Mini-segi-susun-3-Flower-Rak-bumbu-kosmetik-serbaguna-i.145589728.21491056248", "2025-05-25 01:08", 5, "Variasi: Rak segi Flower", "Mohon lebih d tingkatkan lg pelayanannya agar lebih teliti.. Makasih barang2 nya murah² tp berkualitas sayang nya jauh banget klo dr tempt aku mh jd ongkirnya mahal banget walaupun sdh dapet voucer gratis . Sukses terus 👍")
    ,@("(1pcs) Centong entong nasi plastik tebal warna dx termurahVariasi: Ecer_1bj centong PLG", "/(1pcs)-Centong-entong-nasi-plastik-tebal-warna-dx-termurah-i.145589728.2821981447", "2025-05-25 00:17", 1, "Variasi: Ecer_1bj centong PLG", "Ketebalan: penipuDesain: penipuKualitas: penipuada 3 barang gk di kirim .. dan di paket tetep bayar full .. di chat gk di bales`" .. emang sering nipu kali ni toko .. tolong yg lain mending beli di toko oren yg lain aja biar gk ketipu juga")
    ,@("Gelas stand dudukan gelas oval capucino isi 6 gelas plastikVariasi: Gelas Stand Coklat", "/Gelas-stand-dudukan-gelas-oval-capucino-isi-6-gelas-plastik-i.145589728.19090884439", "2025-05-24 23:51", 5, "Variasi: Gelas Stand Coklat", "Desain: bagusKualitas Bahan: okInstalasi: mudahSangat rekomended untuk langganan")
    ,@("Hanger anak plastik benmore 1 pak isi 6pcs grosirVariasi: Hngr benmore coklat", "/Hanger-anak-plastik-benmore-1-pak-isi-6pcs-grosir-i.145589728.28414764389", "2025-05-24 23:23", 5, "Variasi: Hngr benmore coklat", "Barang sesuai pesanan.. dikemas cepat.. pengiriman cepat ,exp lama")
    ,@("Tong sampah segi kicik hitam 5 liter kecil persegi plastik tebalVariasi: TS. Kicik 5L HITAM,Isi 1 pcs", "/Tong-sampah-segi-kicik-hitam-5-liter-kecil-persegi-plastik-tebal-i.145589728.25059427941", "2025-05-24 23:21", 5, "Variasi: TS. Kicik 5L HITAM,Isi 1 pcs", "Desain: TopUkuran: TOPKualitas Bahan: tOPtop")
    ,@("Piring seng enamel lurik ukuran 22 cm besar tebal merk TBVariasi: Piring seng lurik,22cm Besar", "/Piring-seng-enamel-lurik-ukuran-22-cm-besar-tebal-merk-TB-i.145589728.18783036648", "2025-05-24 22:59", 5, "Variasi: Piring seng lurik,22cm Besar", "Ketebalan: sedangDesain: jaduuuullllllKualitas: baikTerima kasih barang sudah sampai dan sesuai")
    ,@("Termos Capucino 6L dan 8L - Rice Ice bucket ekomonis grosirVariasi: Termos eknomis 6L", "/Termos-Capucino-6L-dan-8L-Rice-Ice-bucket-ekomonis-grosir-i.145589728.29500951143", "2025-05-24 22:15", 5, "Variasi: Termos eknomis 6L", "Ukuran: imut imutDesain: bagus warna sesuai picKondisi: dalam kondisi baikTerimakasih seller dan shopee")
    ,@("Kursi jongkok CM Abu-abu warna tebal serbaguna - kursi duduk dapur kokohVariasi: Kursi abu ONIX", "/Kursi-jongkok-CM-Abu-abu-warna-tebal-serbaguna-kursi-duduk-dapur-kokoh-i.145589728.7958001380", "2025-05-24 21:01", 5, "Variasi: Kursi abu ONIX", "Bagus")
    ,@("kapstok baju stainless 6 kait - gantungan baju stainless muranVariasi: Kapstok jempol_6", "/kapstok-baju-stainless-6-kait-gantungan-baju-stainless-muran-i.145589728.6442847024", "2025-05-24 20:37", 5, "Variasi: Kapstok jempol_6", "Desain: eleganDurabilitas: sangat bagusKualitas Bahan: kuat")
    ,@("Botol saos kecap memeclub tinggi 400 ml plastik grosirVariasi: Meme club tinggi", "/Botol-saos-kecap-memeclub-tinggi-400-ml-plastik-grosir-i.145589728.3913320949", "2025-05-24 20:00", 5, "Variasi: Meme club tinggi", "Desain: bagusWarna: sesuaiKualitas: bagus")
    ,@("kapstok baju stainless 6 kait - gantungan baju stainless muranVariasi: Kapstok jempol_6", "/kapstok-baju-stainless-6-kait-gantungan-baju-stainless-muran-i.145589728.6442847024", "2025-05-24 19:59", 5, "Variasi: Kapstok jempol_6", "Desain: bagusDurabilitas: okKualitas Bahan: kuatPuas dengan pesanan saya `nSesuai gambar.,.")
    ,@("(12pcs) Cetakan roti kukus bentuk love rege cerah plastik dxVariasi: 12BJ Ctkn LOVE", "/(12pcs)-Cetakan-roti-kukus-bentuk-love-rege-cerah-plastik-dx-i.145589728.15424268542", "2025-05-24 19:49", 5, "Variasi: 12BJ Ctkn LOVE", "Alhamdulillah paketku sampai dg selamat,barangnya bagus sesuai pict..makasih kak")
    ,@("(12pak) Cotten bud - cotton bud motif bunga halus dan lembut grosirVariasi: 1pak cotenbud bunga", "/(12pak)-Cotten-bud-cotton-bud-motif-bunga-halus-dan-lembut-grosir-i.145589728.16880893501", "2025-05-24 19:45", 5, "Variasi: 1pak cotenbud bunga", "Cocok untuk: semuaKapasitas: bagusKualitas: aman")
    ,@("Mangkok stainless 14 cm 16 cm tebal - baskom stainless serbaguna grosirVariasi: Mangkok 16 cm", "/Mangkok-stainless-14-cm-16-cm-tebal-baskom-stainless-serbaguna-grosir-i.145589728.10650960952", "2025-05-24 19:45", 5, "Variasi: Mangkok 16 cm", "Ukuran: yang di harapkanKualitas Bahan: bagusDurabilitas: kokoh")
    ,@("Tambahan packing kardus paket tambahan agar safetyVariasi: Kardus saftey", "/Tambahan-packing-kardus-paket-tambahan-agar-safety-i.145589728.6615737043", "2025-05-24 19:45", 5, "Variasi: Kardus saftey", "Ketebalan: aman")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Remove rows 22-24 which no longer exist in the updated data
$ws.Range("A22:F24").EntireRow.Delete()

$ws.Range("A1").Select()